# Swap the order of the two "Red ..." symptom paragraphs in the
# alphabetical symptom list so that
#   "Red spots at the back of the roof of the mouth"
# comes immediately before
#   "Red, swollen tonsils"
# (it previously came immediately after it).

$d = $word.ActiveDocument

$textA = "Red, swollen tonsils"
$textB = "Red spots at the back of the roof of the mouth"

# Locate the paragraph containing $textA.
$firstPara = $null
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd("`r", "`a")
    if ($ptext -eq $textA) {
        $firstPara = $p
        break
    }
}

if ($firstPara -ne $null) {
    $secondPara = $firstPara.Next()
    $secondText = $secondPara.Range.Text.TrimEnd("`r", "`a")

    if ($secondText -eq $textB) {
        # Swap the two runs' text, leaving each paragraph's own
        # formatting (rPr / pPr) and paragraph mark untouched.
        $firstPara.Range.Text  = $textB
        $secondPara.Range.Text = $textA
    }
}
